# Update the cached text of the "datetimeFigureOut" date placeholder field
# (Insert > Header & Footer > Date and time) from 8/15/2018 to 12/13/2018
# across the slide master, every slide layout, and the notes master -
# mirroring PowerPoint's behaviour of refreshing every date-placeholder
# field's cached text whenever the deck is saved/republished.

$ppDateTimeNewText = "12/13/2018"
$ppPlaceholderDate = 16

function Update-DatePlaceholderShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePlaceholder = $false
        if ($sh.Type -eq 14) {
            try {
                if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDatePlaceholder = $true
                }
            } catch {
                $isDatePlaceholder = $false
            }
        }
        if ($isDatePlaceholder -and $sh.HasTextFrame) {
            $sh.TextFrame.TextRange.Text = $ppDateTimeNewText
        }
    }
}

$p = $ppt.ActivePresentation

# Slide master
Update-DatePlaceholderShapes $p.SlideMaster.Shapes

# Every slide layout that hangs off the slide master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholderShapes $layouts.Item($li).Shapes
}

# Notes master - its date placeholder shape does not accept direct
# TextFrame writes in this host, but the HeadersFooters facade does.
$p.NotesMaster.HeadersFooters.DateAndTime.Text = $ppDateTimeNewText
